# Scheduled-runner refresh of per-leve marketboard pricing/profit figures.
# Updates currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ)
# columns (H-N) for the affected leve rows on each job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 13216.777
$ws.Range("J17").Value = 13216.777
$ws.Range("L17").Value = 39650.331
$ws.Range("N17").Value = -39986.331

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 12465.333
$ws.Range("I86").Value = 10694.5
$ws.Range("J86").Value = 13350.75
$ws.Range("K86").Value = 10694.5
$ws.Range("L86").Value = 13350.75
$ws.Range("M86").Value = -9571.5
$ws.Range("N86").Value = -15596.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 12465.333
$ws.Range("I89").Value = 10694.5
$ws.Range("J89").Value = 13350.75
$ws.Range("K89").Value = 53472.5
$ws.Range("L89").Value = 66753.75
$ws.Range("M89").Value = -47856.5
$ws.Range("N89").Value = -77985.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 957.9
$ws.Range("I98").Value = 960.125
$ws.Range("J98").Value = 949
$ws.Range("K98").Value = 960.125
$ws.Range("L98").Value = 949
$ws.Range("M98").Value = 537.875
$ws.Range("N98").Value = -3945

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 957.9
$ws.Range("I122").Value = 960.125
$ws.Range("J122").Value = 949
$ws.Range("K122").Value = 2880.375
$ws.Range("L122").Value = 2847
$ws.Range("M122").Value = -430.375
$ws.Range("N122").Value = -7747

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 898.55554
$ws.Range("I137").Value = 898.55554
$ws.Range("K137").Value = 2695.66662
$ws.Range("M137").Value = -145.66662

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3308.5527
$ws.Range("I138").Value = 1691.4166
$ws.Range("J138").Value = 4054.923
$ws.Range("K138").Value = 5074.2498
$ws.Range("L138").Value = 12164.769
$ws.Range("M138").Value = 65.7502000000004
$ws.Range("N138").Value = -22444.769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2569170.5
$ws.Range("I32").Value = 2502647
$ws.Range("K32").Value = 2502647
$ws.Range("M32").Value = -2502360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2767.4211
$ws.Range("I45").Value = 2348.8125
$ws.Range("J45").Value = 5000
$ws.Range("K45").Value = 2348.8125
$ws.Range("L45").Value = 5000
$ws.Range("M45").Value = -1971.8125
$ws.Range("N45").Value = -5754

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 951.6667
$ws.Range("I74").Value = 891.375
$ws.Range("K74").Value = 891.375
$ws.Range("M74").Value = -17.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 951.6667
$ws.Range("I77").Value = 891.375
$ws.Range("K77").Value = 4456.875
$ws.Range("M77").Value = -88.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1020.53845
$ws.Range("I97").Value = 1023.0833
$ws.Range("J97").Value = 990
$ws.Range("K97").Value = 1023.0833
$ws.Range("L97").Value = 990
$ws.Range("M97").Value = -527.0833
$ws.Range("N97").Value = -1982

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1577.9
$ws.Range("I132").Value = 1460.5
$ws.Range("J132").Value = 2047.5
$ws.Range("K132").Value = 4381.5
$ws.Range("L132").Value = 6142.5
$ws.Range("M132").Value = -1851.5
$ws.Range("N132").Value = -11202.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7504
$ws.Range("J20").Value = 5000
$ws.Range("L20").Value = 5000
$ws.Range("N20").Value = -5494

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3707.8823
$ws.Range("I107").Value = 3463.6667
$ws.Range("J107").Value = 4294
$ws.Range("K107").Value = 3463.6667
$ws.Range("L107").Value = 4294
$ws.Range("M107").Value = -1543.6667
$ws.Range("N107").Value = -8134

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1334.6666
$ws.Range("J2").Value = 2000
$ws.Range("L2").Value = 2000
$ws.Range("N2").Value = -2226

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 801.9474
$ws.Range("I31").Value = 789.1
$ws.Range("J31").Value = 816.2222
$ws.Range("K31").Value = 789.1
$ws.Range("L31").Value = 816.2222
$ws.Range("M31").Value = -494.1
$ws.Range("N31").Value = -1406.2222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 801.9474
$ws.Range("I34").Value = 789.1
$ws.Range("J34").Value = 816.2222
$ws.Range("K34").Value = 789.1
$ws.Range("L34").Value = 816.2222
$ws.Range("M34").Value = -587.1
$ws.Range("N34").Value = -1220.2222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H81").Value = 78275
$ws.Range("J81").Value = 78275
$ws.Range("L81").Value = 78275
$ws.Range("N81").Value = -80271

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H84").Value = 78275
$ws.Range("J84").Value = 78275
$ws.Range("L84").Value = 234825
$ws.Range("N84").Value = -244809

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 14956.833
$ws.Range("J88").Value = 14956.833
$ws.Range("L88").Value = 14956.833
$ws.Range("N88").Value = -15768.833

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 14956.833
$ws.Range("J91").Value = 14956.833
$ws.Range("L91").Value = 14956.833
$ws.Range("N91").Value = -17764.833

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 6354.2666
$ws.Range("I132").Value = 6354.2666
$ws.Range("K132").Value = 19062.7998
$ws.Range("M132").Value = -16532.7998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H95").Value = 8878.25
$ws.Range("J95").Value = 8878.25
$ws.Range("L95").Value = 26634.75
$ws.Range("N95").Value = -30752.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 3083.111
$ws.Range("I129").Value = 2929.3333
$ws.Range("J129").Value = 3160
$ws.Range("K129").Value = 8787.999899999999
$ws.Range("L129").Value = 9480
$ws.Range("M129").Value = -3787.999899999999
$ws.Range("N129").Value = -19480

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 42500
$ws.Range("J52").Value = 42500
$ws.Range("L52").Value = 42500
$ws.Range("N52").Value = -43018

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5221.75
$ws.Range("I80").Value = 4995.6665
$ws.Range("K80").Value = 4995.6665
$ws.Range("M80").Value = -3997.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 5221.75
$ws.Range("I83").Value = 4995.6665
$ws.Range("K83").Value = 24978.3325
$ws.Range("M83").Value = -19986.3325

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 28875.25
$ws.Range("J92").Value = 28875.25
$ws.Range("L92").Value = 28875.25
$ws.Range("N92").Value = -32619.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 565.6667
$ws.Range("I97").Value = 298.5
$ws.Range("K97").Value = 298.5
$ws.Range("M97").Value = 197.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3304.125
$ws.Range("I132").Value = 2998.4285
$ws.Range("K132").Value = 8995.2855
$ws.Range("M132").Value = -6465.2855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 150000
$ws.Range("I135").Value = 150000
$ws.Range("J135").Value = 150000
$ws.Range("K135").Value = 150000
$ws.Range("L135").Value = 150000
$ws.Range("M135").Value = -144930
$ws.Range("N135").Value = -160140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6790
$ws.Range("I7").Value = 5697.3335
$ws.Range("J7").Value = 7882.6665
$ws.Range("K7").Value = 5697.3335
$ws.Range("L7").Value = 7882.6665
$ws.Range("M7").Value = -5585.3335
$ws.Range("N7").Value = -8106.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 829.7273
$ws.Range("I22").Value = 811.1429000000001
$ws.Range("J22").Value = 862.25
$ws.Range("K22").Value = 811.1429000000001
$ws.Range("L22").Value = 862.25
$ws.Range("M22").Value = -516.1429000000001
$ws.Range("N22").Value = -1452.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 829.7273
$ws.Range("I27").Value = 811.1429000000001
$ws.Range("J27").Value = 862.25
$ws.Range("K27").Value = 811.1429000000001
$ws.Range("L27").Value = 862.25
$ws.Range("M27").Value = -704.1429000000001
$ws.Range("N27").Value = -1076.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1321.7858
$ws.Range("I55").Value = 1128.125
$ws.Range("K55").Value = 1128.125
$ws.Range("M55").Value = -955.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 993.88
$ws.Range("I93").Value = 1159.1875
$ws.Range("J93").Value = 700
$ws.Range("K93").Value = 1159.1875
$ws.Range("L93").Value = 700
$ws.Range("M93").Value = 88.8125
$ws.Range("N93").Value = -3196

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 6790
$ws.Range("I126").Value = 5697.3335
$ws.Range("J126").Value = 7882.6665
$ws.Range("K126").Value = 17092.0005
$ws.Range("L126").Value = 23647.9995
$ws.Range("M126").Value = -14622.0005
$ws.Range("N126").Value = -28587.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2820.9048
$ws.Range("I132").Value = 2740.2856
$ws.Range("J132").Value = 2982.1428
$ws.Range("K132").Value = 8220.856800000001
$ws.Range("L132").Value = 8946.428400000001
$ws.Range("M132").Value = -5690.856800000001
$ws.Range("N132").Value = -14006.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1631.6
$ws.Range("J136").Value = 1498.5
$ws.Range("L136").Value = 4495.5
$ws.Range("N136").Value = -9595.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 28749.5
$ws.Range("I52").Value = 20000
$ws.Range("J52").Value = 31666
$ws.Range("K52").Value = 20000
$ws.Range("L52").Value = 31666
$ws.Range("M52").Value = -19774
$ws.Range("N52").Value = -32118

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2136.923
$ws.Range("I126").Value = 1744
$ws.Range("J126").Value = 2311.5557
$ws.Range("K126").Value = 5232
$ws.Range("L126").Value = 6934.6671
$ws.Range("M126").Value = -2762
$ws.Range("N126").Value = -11874.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2165.3333
$ws.Range("I132").Value = 2748
$ws.Range("K132").Value = 8244
$ws.Range("M132").Value = -5714
